# Update "想去人数" (want-to-go count) values in column F across all sheets
# as published in the refreshed gh-pages data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1220
$ws1.Range("F4").Value = 49
$ws1.Range("F5").Value = 1391
$ws1.Range("F6").Value = 1722
$ws1.Range("F7").Value = 6248
$ws1.Range("F9").Value = 1832
$ws1.Range("F10").Value = 488
$ws1.Range("F12").Value = 21
$ws1.Range("F16").Value = 7013
$ws1.Range("F17").Value = 129
$ws1.Range("F21").Value = 1715
$ws1.Range("F22").Value = 842
$ws1.Range("F23").Value = 20
$ws1.Range("F25").Value = 166
$ws1.Range("F26").Value = 1620
$ws1.Range("F27").Value = 770
$ws1.Range("F28").Value = 326
$ws1.Range("F31").Value = 60
$ws1.Range("F32").Value = 84

# Sheet 2
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 351
$ws2.Range("F5").Value = 202

# Sheet 3
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 9533

# Sheet 4
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 9533
$ws4.Range("F5").Value = 1220
$ws4.Range("F7").Value = 49
$ws4.Range("F9").Value = 351
$ws4.Range("F10").Value = 1392
$ws4.Range("F12").Value = 1722
$ws4.Range("F13").Value = 6248
$ws4.Range("F14").Value = 1832
$ws4.Range("F17").Value = 488
$ws4.Range("F19").Value = 21
$ws4.Range("F24").Value = 7013
$ws4.Range("F25").Value = 129
$ws4.Range("F29").Value = 1715
$ws4.Range("F30").Value = 842
$ws4.Range("F31").Value = 20
$ws4.Range("F32").Value = 166
$ws4.Range("F33").Value = 1620
$ws4.Range("F34").Value = 770
$ws4.Range("F36").Value = 326
$ws4.Range("F42").Value = 84
